$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4817
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 6156
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 6156
$ws.Range("M2").Value = -687
$ws.Range("N2").Value = -6382
$ws.Range("H62").Value = 7799.6
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376
$ws.Range("H65").Value = 7799.6
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880
$ws.Range("H86").Value = 5871.25
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 6138.5713
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 6138.5713
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -8384.5713
$ws.Range("H89").Value = 5871.25
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 6138.5713
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 30692.8565
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -41924.85649999999
$ws.Range("H113").Value = 3612.25
$ws.Range("I113").Value = 3399.6667
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 3399.6667
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = -145.6667000000002
$ws.Range("N113").Value = -10758
$ws.Range("H138").Value = 3902.7542
$ws.Range("J138").Value = 4008.6667
$ws.Range("L138").Value = 12026.0001
$ws.Range("N138").Value = -22306.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8729
$ws.Range("I32").Value = 6568.0454
$ws.Range("J32").Value = 32499.5
$ws.Range("K32").Value = 6568.0454
$ws.Range("L32").Value = 32499.5
$ws.Range("M32").Value = -6281.0454
$ws.Range("N32").Value = -33073.5
$ws.Range("H45").Value = 1107.5555
$ws.Range("I45").Value = 1196
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 1196
$ws.Range("L45").Value = 400
$ws.Range("M45").Value = -819
$ws.Range("N45").Value = -1154
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240
$ws.Range("H110").Value = 3000
$ws.Range("J110").Value = 3000
$ws.Range("L110").Value = 3000
$ws.Range("N110").Value = -7090
$ws.Range("H122").Value = 2509.8572
$ws.Range("J122").Value = 2709.6
$ws.Range("L122").Value = 8128.799999999999
$ws.Range("N122").Value = -13028.8
$ws.Range("H127").Value = 127499
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 127499
$ws.Range("K127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("M127").Value = 127499
$ws.Range("N127").Value = -137419

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 33996
$ws.Range("J6").Value = 33996
$ws.Range("L6").Value = 33996
$ws.Range("N6").Value = -34222
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0
$ws.Range("H86").Value = 2016.6
$ws.Range("I86").Value = 1959.75
$ws.Range("J86").Value = 2244
$ws.Range("K86").Value = 1959.75
$ws.Range("L86").Value = 2244
$ws.Range("M86").Value = -836.75
$ws.Range("N86").Value = -4490
$ws.Range("H89").Value = 2016.6
$ws.Range("I89").Value = 1959.75
$ws.Range("J89").Value = 2244
$ws.Range("K89").Value = 9798.75
$ws.Range("L89").Value = 11220
$ws.Range("M89").Value = -4182.75
$ws.Range("N89").Value = -22452
$ws.Range("H105").Value = 2679.5386
$ws.Range("I105").Value = 2736.25
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 2736.25
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -989.25
$ws.Range("N105").Value = -5493
$ws.Range("H134").Value = 899.5
$ws.Range("I134").Value = 899
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 2697
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -162
$ws.Range("N134").Value = -7770

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2191
$ws.Range("J122").Value = 2345.8
$ws.Range("L122").Value = 7037.400000000001
$ws.Range("N122").Value = -11937.4
$ws.Range("H132").Value = 2234.4167
$ws.Range("I132").Value = 1226.875
$ws.Range("K132").Value = 3680.625
$ws.Range("M132").Value = -1150.625
$ws.Range("H134").Value = 3541.077
$ws.Range("I134").Value = 3419.5
$ws.Range("K134").Value = 10258.5
$ws.Range("M134").Value = -7723.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 13580.9
$ws.Range("I94").Value = 4327.25
$ws.Range("K94").Value = 12981.75
$ws.Range("M94").Value = -12305.75
$ws.Range("H113").Value = 1244.3572
$ws.Range("J113").Value = 1059.7273
$ws.Range("L113").Value = 3179.1819
$ws.Range("N113").Value = -7519.1819
$ws.Range("H137").Value = 3857.923
$ws.Range("I137").Value = 2023.4
$ws.Range("J137").Value = 5004.5
$ws.Range("K137").Value = 6070.200000000001
$ws.Range("L137").Value = 15013.5
$ws.Range("M137").Value = -970.2000000000007
$ws.Range("N137").Value = -25213.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 73.63636
$ws.Range("I2").Value = 103.42857
$ws.Range("J2").Value = 21.5
$ws.Range("K2").Value = 103.42857
$ws.Range("L2").Value = 21.5
$ws.Range("M2").Value = 9.571430000000007
$ws.Range("N2").Value = -247.5
$ws.Range("H80").Value = 3443.2
$ws.Range("I80").Value = 2129.6667
$ws.Range("J80").Value = 3771.5833
$ws.Range("K80").Value = 2129.6667
$ws.Range("L80").Value = 3771.5833
$ws.Range("M80").Value = -1131.6667
$ws.Range("N80").Value = -5767.5833
$ws.Range("H83").Value = 3443.2
$ws.Range("I83").Value = 2129.6667
$ws.Range("J83").Value = 3771.5833
$ws.Range("K83").Value = 10648.3335
$ws.Range("L83").Value = 18857.9165
$ws.Range("M83").Value = -5656.333500000001
$ws.Range("N83").Value = -28841.9165
$ws.Range("H113").Value = 1816.6666
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 1951.6957
$ws.Range("I122").Value = 1148.9375
$ws.Range("J122").Value = 3786.5715
$ws.Range("K122").Value = 3446.8125
$ws.Range("L122").Value = 11359.7145
$ws.Range("M122").Value = -996.8125
$ws.Range("N122").Value = -16259.7145
$ws.Range("H132").Value = 1657.5714
$ws.Range("I132").Value = 1577.2354
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 4731.706200000001
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -2201.706200000001
$ws.Range("N132").Value = -11057
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 682.55554
$ws.Range("J22").Value = 561.625
$ws.Range("L22").Value = 561.625
$ws.Range("N22").Value = -1151.625
$ws.Range("H27").Value = 682.55554
$ws.Range("J27").Value = 561.625
$ws.Range("L27").Value = 561.625
$ws.Range("N27").Value = -775.625
$ws.Range("H55").Value = 395
$ws.Range("I55").Value = 395
$ws.Range("K55").Value = 395
$ws.Range("M55").Value = -222
$ws.Range("H82").Value = 2135.9285
$ws.Range("I82").Value = 1768.5555
$ws.Range("J82").Value = 2797.2
$ws.Range("K82").Value = 1768.5555
$ws.Range("L82").Value = 2797.2
$ws.Range("M82").Value = -1407.5555
$ws.Range("N82").Value = -3519.2
$ws.Range("H85").Value = 2135.9285
$ws.Range("I85").Value = 1768.5555
$ws.Range("J85").Value = 2797.2
$ws.Range("K85").Value = 1768.5555
$ws.Range("L85").Value = 2797.2
$ws.Range("M85").Value = -520.5554999999999
$ws.Range("N85").Value = -5293.2
$ws.Range("H122").Value = 6199.9165
$ws.Range("I122").Value = 6108.16
$ws.Range("J122").Value = 6408.4546
$ws.Range("K122").Value = 18324.48
$ws.Range("L122").Value = 19225.3638
$ws.Range("M122").Value = -15874.48
$ws.Range("N122").Value = -24125.3638
$ws.Range("H132").Value = 1399
$ws.Range("I132").Value = 1399
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4197
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 834.6667
$ws.Range("I17").Value = 252
$ws.Range("K17").Value = 252
$ws.Range("M17").Value = -80
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H98").Value = 48000
$ws.Range("J98").Value = 48000
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -53990
$ws.Range("H107").Value = 1038.4375
$ws.Range("J107").Value = 1041.3334
$ws.Range("L107").Value = 3124.0002
$ws.Range("N107").Value = -6964.0002
$ws.Range("H126").Value = 5204.9287
$ws.Range("I126").Value = 3456.6
$ws.Range("J126").Value = 6176.222
$ws.Range("K126").Value = 10369.8
$ws.Range("L126").Value = 18528.666
$ws.Range("M126").Value = -7899.799999999999
$ws.Range("N126").Value = -23468.666
$ws.Range("H130").Value = 56499
$ws.Range("J130").Value = 56499
$ws.Range("L130").Value = 56499
$ws.Range("N130").Value = -66539
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0
